# "feat - ajustado tela de acesso"
#
# The "@prints" placeholder moves from the "Problema" table into the
# "Solucao" table (replacing "@image1" there), merging into a single,
# formatted paragraph; the run left behind in the "Problema" table is
# removed while keeping its (now empty, but still formatted) paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Locate the paragraph that currently holds "@image1" (Solucao
#    table) and the blank, already-formatted paragraph right after it.
# ---------------------------------------------------------------------
$imagePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "@image1") {
        $imagePara = $p
        break
    }
}

if ($null -eq $imagePara) {
    throw "Could not find @image1 placeholder"
}

$nextPara = $imagePara.Next()

# Replace the "@image1" paragraph with a single paragraph that carries
# the formatting (Arial, bold, sz20) previously living on the blank
# paragraph that followed it, and the new "@prints" text.
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="SemEspaamento"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:t>@prints</w:t></w:r></w:p>'
[void]$imagePara.Range.InsertXML($newParaXml)

# Remove the now-redundant blank paragraph that used to carry the
# formatting (its rPr has been folded into the paragraph above).
[void]$nextPara.Range.Delete()

# ---------------------------------------------------------------------
# 2) Locate the remaining "@prints" run (Problema table) and delete
#    just the run's text, leaving the (now empty) paragraph and all of
#    its paragraph-mark formatting untouched.
# ---------------------------------------------------------------------
$printsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "@prints" -and $p.Range.Start -ne $imagePara.Range.Start) {
        $printsPara = $p
        break
    }
}

if ($null -eq $printsPara) {
    throw "Could not find @prints placeholder"
}

$textOnly = $d.Range($printsPara.Range.Start, $printsPara.Range.End - 1)
[void]$textOnly.Delete()
